{"js": "// The document contains a single 20x5 table of arithmetic expressions\n// (e.g. \"71-19=\"). The commit replaces each cell's expression with a new\n// one, in document (row-major) order. Duplicate \"before\" values exist\n// (e.g. \"0+0=\" appears twice, \"29-22=\" appears twice) with DIFFERENT\n// replacements, so the update must be positional (row/column index),\n// not a global text find-and-replace.\nconst newValues = [\n  [\"68+15=\", \"13+77=\", \"31-13=\", \"23-9=\", \"47-28=\"],\n  [\"76-73=\", \"61+1=\", \"74-0=\", \"31-1=\", \"36-9=\"],\n  [\"45-16=\", \"92-49=\", \"9+80=\", \"84-55=\", \"10+26=\"],\n  [\"53-11=\", \"64-49=\", \"89-71=\", \"44-22=\", \"35+59=\"],\n  [\"91-30=\", \"17-17=\", \"71+26=\", \"22+33=\", \"20+79=\"],\n  [\"26-22=\", \"45+37=\", \"10+13=\", \"70-54=\", \"17+51=\"],\n  [\"35+5=\", \"97-8=\", \"27+58=\", \"72+0=\", \"10+0=\"],\n  [\"83-46=\", \"77-13=\", \"35+34=\", \"20+77=\", \"57+16=\"],\n  [\"14+15=\", \"76+11=\", \"35+4=\", \"68-28=\", \"13+38=\"],\n  [\"13-8=\", \"16+13=\", \"96-3=\", \"48-1=\", \"77-35=\"],\n  [\"83+15=\", \"35+60=\", \"37+59=\", \"81-28=\", \"2+76=\"],\n  [\"52-51=\", \"46+48=\", \"18-10=\", \"55+30=\", \"78+0=\"],\n  [\"18+15=\", \"61-59=\", \"79-50=\", \"67+28=\", \"47-9=\"],\n  [\"62+28=\", \"90-69=\", \"68+18=\", \"26+33=\", \"35-22=\"],\n  [\"40-7=\", \"89-78=\", \"66-2=\", \"83-35=\", \"56-17=\"],\n  [\"71+7=\", \"39+10=\", \"15+6=\", \"85+2=\", \"10+49=\"],\n  [\"83-19=\", \"99-40=\", \"57+4=\", \"86-20=\", \"28-0=\"],\n  [\"88-23=\", \"64-8=\", \"27+71=\", \"58-5=\", \"89-72=\"],\n  [\"99-75=\", \"9+51=\", \"99-51=\", \"8+47=\", \"24-13=\"],\n  [\"31+4=\", \"89-29=\", \"62-15=\", \"87-46=\", \"28+1=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n// Assigning `.values` rewrites each cell's text in place while leaving\n// the existing run/paragraph formatting (font, size, etc.) untouched.\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# The document contains a single 20x5 table of arithmetic expressions\n# (e.g. \"71-19=\"). The commit replaces each cell's expression with a new\n# one, in document (row/column) order. Duplicate \"before\" values exist\n# (e.g. \"0+0=\" appears twice, \"29-22=\" appears twice) with DIFFERENT\n# replacements, so the update must be positional (row/column index),\n# not a global text find-and-replace.\n$newValues = @(\n    @(\"68+15=\", \"13+77=\", \"31-13=\", \"23-9=\", \"47-28=\"),\n    @(\"76-73=\", \"61+1=\", \"74-0=\", \"31-1=\", \"36-9=\"),\n    @(\"45-16=\", \"92-49=\", \"9+80=\", \"84-55=\", \"10+26=\"),\n    @(\"53-11=\", \"64-49=\", \"89-71=\", \"44-22=\", \"35+59=\"),\n    @(\"91-30=\", \"17-17=\", \"71+26=\", \"22+33=\", \"20+79=\"),\n    @(\"26-22=\", \"45+37=\", \"10+13=\", \"70-54=\", \"17+51=\"),\n    @(\"35+5=\", \"97-8=\", \"27+58=\", \"72+0=\", \"10+0=\"),\n    @(\"83-46=\", \"77-13=\", \"35+34=\", \"20+77=\", \"57+16=\"),\n    @(\"14+15=\", \"76+11=\", \"35+4=\", \"68-28=\", \"13+38=\"),\n    @(\"13-8=\", \"16+13=\", \"96-3=\", \"48-1=\", \"77-35=\"),\n    @(\"83+15=\", \"35+60=\", \"37+59=\", \"81-28=\", \"2+76=\"),\n    @(\"52-51=\", \"46+48=\", \"18-10=\", \"55+30=\", \"78+0=\"),\n    @(\"18+15=\", \"61-59=\", \"79-50=\", \"67+28=\", \"47-9=\"),\n    @(\"62+28=\", \"90-69=\", \"68+18=\", \"26+33=\", \"35-22=\"),\n    @(\"40-7=\", \"89-78=\", \"66-2=\", \"83-35=\", \"56-17=\"),\n    @(\"71+7=\", \"39+10=\", \"15+6=\", \"85+2=\", \"10+49=\"),\n    @(\"83-19=\", \"99-40=\", \"57+4=\", \"86-20=\", \"28-0=\"),\n    @(\"88-23=\", \"64-8=\", \"27+71=\", \"58-5=\", \"89-72=\"),\n    @(\"99-75=\", \"9+51=\", \"99-51=\", \"8+47=\", \"24-13=\"),\n    @(\"31+4=\", \"89-29=\", \"62-15=\", \"87-46=\", \"28+1=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
